$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to match repulled data / mean calculation
$ws.Range("F2").Value = -6
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -4
$ws.Range("F13").Value = 0
